$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.873.58'
$ws.Range('E2').Value = '  -5.24%  '

$ws.Range('D3').Value = '2.524.93'
$ws.Range('E3').Value = '  -4.75%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '502.69'
$ws.Range('E5').Value = '  -5.74%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.03'
$ws.Range('E6').Value = '  -8.34%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.25%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.560'
$ws.Range('E8').Value = '  -5.34%  '

$ws.Range('D9').Value = '2.519.12'
$ws.Range('E9').Value = '  -5.42%  '

$ws.Range('E10').Value = '  -8.10%  '

$ws.Range('E11').Value = '  -7.90%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.328'
$ws.Range('E12').Value = '  -6.76%  '

$ws.Range('E13').Value = '  -0.72%  '

$ws.Range('D14').Value = '2.971.88'
$ws.Range('E14').Value = '  -4.56%  '

$ws.Range('D15').Value = '57.911.57'
$ws.Range('E15').Value = '  -5.15%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.55'
$ws.Range('E16').Value = '  -7.00%  '

$ws.Range('E17').Value = '  -7.33%  '

$ws.Range('D18').Value = '2.533.92'
$ws.Range('E18').Value = '  -4.75%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.47'
$ws.Range('E19').Value = '  -6.39%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.54'
$ws.Range('E20').Value = '  -4.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.98'
$ws.Range('E21').Value = '  -6.65%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.88'
$ws.Range('E23').Value = '  -5.83%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.18'
$ws.Range('E24').Value = '  -2.48%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.406'
$ws.Range('E25').Value = '  -5.89%  '

$ws.Range('E26').Value = '  +0.18%  '

$ws.Range('D27').Value = '2.650.79'
$ws.Range('E27').Value = '  -4.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.156'
$ws.Range('E28').Value = '  -7.27%  '

$ws.Range('E29').Value = '  -9.89%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.87'
$ws.Range('E30').Value = '  -7.05%  '

$ws.Range('E31').Value = '  +0.15%  '

$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.84'
$ws.Range('E32').Value = '  -5.88%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '148.93'
$ws.Range('E33').Value = '  -0.84%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.34'
$ws.Range('E34').Value = '  -6.21%  '

$ws.Range('E35').Value = '  -7.20%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.85'
$ws.Range('E36').Value = '  -7.08%  '

$ws.Range('E37').Value = '  -0.61%  '

$ws.Range('E38').Value = '  -8.78%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.71'
$ws.Range('E39').Value = '  -3.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.807'
$ws.Range('E40').Value = '  -12.80%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.36'
$ws.Range('E41').Value = '  -9.23%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.48'
$ws.Range('E42').Value = '  -8.97%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('E43').Value = '  -0.17%  '

$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '276.71'
$ws.Range('E44').Value = '  -10.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0985'
$ws.Range('E45').Value = '  -3.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.591'
$ws.Range('E46').Value = '  -8.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0527'
$ws.Range('E47').Value = '  -6.67%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.27'
$ws.Range('E48').Value = '  -0.78%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.44'
$ws.Range('E49').Value = '  -7.73%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0223'
$ws.Range('E50').Value = '  -6.59%  '

$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.888.63'
$ws.Range('E51').Value = '  -5.37%  '
